$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2892.375
$ws.Range("I19").Value = 4013.8333
$ws.Range("K19").Value = 4013.8333
$ws.Range("M19").Value = -3838.8333
$ws.Range("H40").Value = 2868.6924
$ws.Range("I40").Value = 4119.5293
$ws.Range("J40").Value = 1902.1364
$ws.Range("K40").Value = 4119.5293
$ws.Range("L40").Value = 1902.1364
$ws.Range("M40").Value = -3944.5293
$ws.Range("N40").Value = -2252.1364
$ws.Range("H53").Value = 362.92307
$ws.Range("I53").Value = 358.33334
$ws.Range("K53").Value = 358.33334
$ws.Range("M53").Value = 278.66666
$ws.Range("H55").Value = 300.5238
$ws.Range("I55").Value = 201.90909
$ws.Range("K55").Value = 201.90909
$ws.Range("M55").Value = 12.09091000000001
$ws.Range("H70").Value = 2178.5557
$ws.Range("I70").Value = 1628.2222
$ws.Range("K70").Value = 4884.6666
$ws.Range("M70").Value = -4614.6666
$ws.Range("H73").Value = 2178.5557
$ws.Range("I73").Value = 1628.2222
$ws.Range("K73").Value = 4884.6666
$ws.Range("M73").Value = -3948.6666
$ws.Range("H92").Value = 804.88
$ws.Range("I92").Value = 425.33334
$ws.Range("J92").Value = 2797.5
$ws.Range("K92").Value = 425.33334
$ws.Range("L92").Value = 2797.5
$ws.Range("M92").Value = 822.66666
$ws.Range("N92").Value = -5293.5
$ws.Range("H116").Value = 199300.77
$ws.Range("I116").Value = 46537.152
$ws.Range("J116").Value = 447541.62
$ws.Range("K116").Value = 46537.152
$ws.Range("L116").Value = 447541.62
$ws.Range("M116").Value = -43095.152
$ws.Range("N116").Value = -454425.62
$ws.Range("H132").Value = 72554.82000000001
$ws.Range("I132").Value = 76751.91
$ws.Range("K132").Value = 230255.73
$ws.Range("M132").Value = -227725.73
$ws.Range("H138").Value = 2465.0732
$ws.Range("I138").Value = 1145.4286
$ws.Range("J138").Value = 5307.385
$ws.Range("K138").Value = 3436.2858
$ws.Range("L138").Value = 15922.155
$ws.Range("M138").Value = 1703.7142
$ws.Range("N138").Value = -26202.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 207.5
$ws.Range("I4").Value = 208.18182
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 208.18182
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -92.18181999999999
$ws.Range("N4").Value = -432
$ws.Range("H32").Value = 21129.836
$ws.Range("I32").Value = 22250.79
$ws.Range("K32").Value = 22250.79
$ws.Range("M32").Value = -21963.79
$ws.Range("H61").Value = 795942.6
$ws.Range("I61").Value = 856966.9
$ws.Range("K61").Value = 856966.9
$ws.Range("M61").Value = -856754.9
$ws.Range("H97").Value = 871.55554
$ws.Range("I97").Value = 896.5294
$ws.Range("J97").Value = 447
$ws.Range("K97").Value = 896.5294
$ws.Range("L97").Value = 447
$ws.Range("M97").Value = -400.5294
$ws.Range("N97").Value = -1439
$ws.Range("H110").Value = 1824.75
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1824.75
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 1824.75
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -5914.75
$ws.Range("H136").Value = 795942.6
$ws.Range("I136").Value = 856966.9
$ws.Range("K136").Value = 2570900.7
$ws.Range("M136").Value = -2568350.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1005.5
$ws.Range("I20").Value = 1007.4667
$ws.Range("J20").Value = 995.6667
$ws.Range("K20").Value = 1007.4667
$ws.Range("L20").Value = 995.6667
$ws.Range("M20").Value = -760.4666999999999
$ws.Range("N20").Value = -1489.6667
$ws.Range("H35").Value = 17000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H37").Value = 126
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 954.73914
$ws.Range("I22").Value = 713.5
$ws.Range("J22").Value = 1506.1428
$ws.Range("K22").Value = 713.5
$ws.Range("L22").Value = 1506.1428
$ws.Range("M22").Value = -363.5
$ws.Range("N22").Value = -2206.1428
$ws.Range("H51").Value = 32855.43
$ws.Range("J51").Value = 82494
$ws.Range("L51").Value = 82494
$ws.Range("N51").Value = -83966
$ws.Range("H58").Value = 540108.4
$ws.Range("I58").Value = 774198.9
$ws.Range("K58").Value = 774198.9
$ws.Range("M58").Value = -773995.9
$ws.Range("H61").Value = 32855.43
$ws.Range("J61").Value = 82494
$ws.Range("L61").Value = 82494
$ws.Range("N61").Value = -83190
$ws.Range("H74").Value = 50104.332
$ws.Range("J74").Value = 50104.332
$ws.Range("L74").Value = 50104.332
$ws.Range("N74").Value = -51852.332
$ws.Range("H77").Value = 50104.332
$ws.Range("J77").Value = 50104.332
$ws.Range("L77").Value = 150312.996
$ws.Range("N77").Value = -159048.996
$ws.Range("H122").Value = 3775.6
$ws.Range("I122").Value = 2251.3333
$ws.Range("J122").Value = 5022.727
$ws.Range("K122").Value = 6753.999899999999
$ws.Range("L122").Value = 15068.181
$ws.Range("M122").Value = -4303.999899999999
$ws.Range("N122").Value = -19968.181
$ws.Range("H123").Value = 129088.55
$ws.Range("J123").Value = 129088.55
$ws.Range("L123").Value = 129088.55
$ws.Range("N123").Value = -138888.55
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H134").Value = 2024.1666
$ws.Range("I134").Value = 1849.1
$ws.Range("K134").Value = 5547.299999999999
$ws.Range("M134").Value = -3012.299999999999
$ws.Range("H136").Value = 540108.4
$ws.Range("I136").Value = 774198.9
$ws.Range("K136").Value = 2322596.7
$ws.Range("M136").Value = -2320046.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 844.3333
$ws.Range("I5").Value = 320.8
$ws.Range("J5").Value = 1498.75
$ws.Range("K5").Value = 962.4000000000001
$ws.Range("L5").Value = 4496.25
$ws.Range("M5").Value = -850.4000000000001
$ws.Range("N5").Value = -4720.25
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H51").Value = 4931.4287
$ws.Range("I51").Value = 2750
$ws.Range("J51").Value = 5804
$ws.Range("K51").Value = 8250
$ws.Range("L51").Value = 17412
$ws.Range("M51").Value = -7790
$ws.Range("N51").Value = -18332
$ws.Range("H131").Value = 11758.762
$ws.Range("I131").Value = 777.0909
$ws.Range("J131").Value = 23838.6
$ws.Range("K131").Value = 2331.2727
$ws.Range("L131").Value = 71515.79999999999
$ws.Range("M131").Value = 2708.7273
$ws.Range("N131").Value = -81595.79999999999
$ws.Range("H135").Value = 844.3333
$ws.Range("I135").Value = 320.8
$ws.Range("J135").Value = 1498.75
$ws.Range("K135").Value = 2887.2
$ws.Range("L135").Value = 13488.75
$ws.Range("M135").Value = -352.2000000000003
$ws.Range("N135").Value = -18558.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1566.2
$ws.Range("I31").Value = 1566.2
$ws.Range("K31").Value = 1566.2
$ws.Range("M31").Value = -1274.2
$ws.Range("H37").Value = 1566.2
$ws.Range("I37").Value = 1566.2
$ws.Range("K37").Value = 1566.2
$ws.Range("M37").Value = -1289.2
$ws.Range("H70").Value = 5663
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 6494.5
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 6494.5
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -7034.5
$ws.Range("H73").Value = 5663
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 6494.5
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 6494.5
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -8366.5
$ws.Range("H107").Value = 21026.824
$ws.Range("I107").Value = 50555.285
$ws.Range("K107").Value = 50555.285
$ws.Range("M107").Value = -48635.285
$ws.Range("H122").Value = 127250
$ws.Range("I122").Value = 127250
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 381750
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -379300
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2514.739
$ws.Range("I40").Value = 2524.5454
$ws.Range("K40").Value = 2524.5454
$ws.Range("M40").Value = -2388.5454
$ws.Range("H130").Value = 88367.57000000001
$ws.Range("J130").Value = 88367.57000000001
$ws.Range("L130").Value = 88367.57000000001
$ws.Range("N130").Value = -98407.57000000001
$ws.Range("H136").Value = 3028.739
$ws.Range("I136").Value = 2127.5293
$ws.Range("K136").Value = 6382.5879
$ws.Range("M136").Value = -3832.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1545.7142
$ws.Range("I81").Value = 1470.1666
$ws.Range("J81").Value = 1999
$ws.Range("K81").Value = 2940.3332
$ws.Range("L81").Value = 3998
$ws.Range("M81").Value = -1879.3332
$ws.Range("N81").Value = -6120
$ws.Range("H84").Value = 1545.7142
$ws.Range("I84").Value = 1470.1666
$ws.Range("J84").Value = 1999
$ws.Range("K84").Value = 14701.666
$ws.Range("L84").Value = 19990
$ws.Range("M84").Value = -9397.666000000001
$ws.Range("N84").Value = -30598
$ws.Range("H96").Value = 3332.6667
$ws.Range("I96").Value = 1999.5
$ws.Range("J96").Value = 5999
$ws.Range("K96").Value = 1999.5
$ws.Range("L96").Value = 5999
$ws.Range("M96").Value = -626.5
$ws.Range("N96").Value = -8745
